$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1111.7273
$ws.Cells.Item(28, 9).Value = 1122.9
$ws.Cells.Item(28, 10).Value = 1000
$ws.Cells.Item(28, 11).Value = 1122.9
$ws.Cells.Item(28, 12).Value = 1000
$ws.Cells.Item(28, 13).Value = -637.9000000000001
$ws.Cells.Item(28, 14).Value = -1970

$ws.Cells.Item(64, 8).Value = 1433.6666
$ws.Cells.Item(64, 10).Value = 1433.6666
$ws.Cells.Item(64, 12).Value = 1433.6666
$ws.Cells.Item(64, 14).Value = -1929.6666

$ws.Cells.Item(67, 8).Value = 1433.6666
$ws.Cells.Item(67, 10).Value = 1433.6666
$ws.Cells.Item(67, 12).Value = 1433.6666
$ws.Cells.Item(67, 14).Value = -3149.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 970.5
$ws.Cells.Item(32, 9).Value = 713.3333
$ws.Cells.Item(32, 11).Value = 713.3333
$ws.Cells.Item(32, 13).Value = -426.3333

$ws.Cells.Item(61, 8).Value = 1729.5
$ws.Cells.Item(61, 9).Value = 1624.7
$ws.Cells.Item(61, 10).Value = 2253.5
$ws.Cells.Item(61, 11).Value = 1624.7
$ws.Cells.Item(61, 12).Value = 2253.5
$ws.Cells.Item(61, 13).Value = -1412.7
$ws.Cells.Item(61, 14).Value = -2677.5

$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 1729.5
$ws.Cells.Item(136, 9).Value = 1624.7
$ws.Cells.Item(136, 10).Value = 2253.5
$ws.Cells.Item(136, 11).Value = 4874.1
$ws.Cells.Item(136, 12).Value = 6760.5
$ws.Cells.Item(136, 13).Value = -2324.1
$ws.Cells.Item(136, 14).Value = -11860.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2000
$ws.Cells.Item(134, 9).Value = 2000
$ws.Cells.Item(134, 11).Value = 6000
$ws.Cells.Item(134, 13).Value = -3465

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9852.375
$ws.Cells.Item(31, 10).Value = 13171
$ws.Cells.Item(31, 12).Value = 13171
$ws.Cells.Item(31, 14).Value = -13761

$ws.Cells.Item(34, 8).Value = 9852.375
$ws.Cells.Item(34, 10).Value = 13171
$ws.Cells.Item(34, 12).Value = 13171
$ws.Cells.Item(34, 14).Value = -13575

$ws.Cells.Item(36, 8).Value = 4999.25
$ws.Cells.Item(36, 9).Value = 4999.25
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 4999.25
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -4611.25
$ws.Cells.Item(36, 14).ClearContents()

$ws.Cells.Item(40, 8).Value = 4999.25
$ws.Cells.Item(40, 9).Value = 4999.25
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 4999.25
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -4839.25
$ws.Cells.Item(40, 14).ClearContents()

$ws.Cells.Item(44, 8).Value = 29999.2
$ws.Cells.Item(44, 9).Value = 29998.5
$ws.Cells.Item(44, 10).Value = 29999.666
$ws.Cells.Item(44, 11).Value = 29998.5
$ws.Cells.Item(44, 12).Value = 29999.666
$ws.Cells.Item(44, 13).Value = -29556.5
$ws.Cells.Item(44, 14).Value = -30883.666

$ws.Cells.Item(58, 8).Value = 1640.4166
$ws.Cells.Item(58, 10).Value = 1637.2
$ws.Cells.Item(58, 12).Value = 1637.2
$ws.Cells.Item(58, 14).Value = -2043.2

$ws.Cells.Item(59, 8).Value = 65000
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 65000
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 12).Value = 65000
$ws.Cells.Item(59, 13).ClearContents()
$ws.Cells.Item(59, 14).Value = -67290

$ws.Cells.Item(132, 8).Value = 3000
$ws.Cells.Item(132, 9).Value = 2000
$ws.Cells.Item(132, 10).Value = 4000
$ws.Cells.Item(132, 11).Value = 6000
$ws.Cells.Item(132, 12).Value = 12000
$ws.Cells.Item(132, 13).Value = -3470
$ws.Cells.Item(132, 14).Value = -17060

$ws.Cells.Item(136, 8).Value = 1640.4166
$ws.Cells.Item(136, 10).Value = 1637.2
$ws.Cells.Item(136, 12).Value = 4911.6
$ws.Cells.Item(136, 14).Value = -10011.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 2828.5715
$ws.Cells.Item(86, 9).Value = 2760
$ws.Cells.Item(86, 10).Value = 3000
$ws.Cells.Item(86, 11).Value = 8280
$ws.Cells.Item(86, 12).Value = 9000
$ws.Cells.Item(86, 13).Value = -7094
$ws.Cells.Item(86, 14).Value = -11372

$ws.Cells.Item(89, 8).Value = 2828.5715
$ws.Cells.Item(89, 9).Value = 2760
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 24840
$ws.Cells.Item(89, 12).Value = 27000
$ws.Cells.Item(89, 13).Value = -18912
$ws.Cells.Item(89, 14).Value = -38856

$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 13).ClearContents()

$ws.Cells.Item(122, 8).Value = 104
$ws.Cells.Item(122, 9).Value = 104
$ws.Cells.Item(122, 11).Value = 936
$ws.Cells.Item(122, 13).Value = 1514

$ws.Cells.Item(131, 8).Value = 4093.5454
$ws.Cells.Item(131, 9).Value = 2405.8
$ws.Cells.Item(131, 10).Value = 5500
$ws.Cells.Item(131, 11).Value = 7217.400000000001
$ws.Cells.Item(131, 12).Value = 16500
$ws.Cells.Item(131, 13).Value = -2177.400000000001
$ws.Cells.Item(131, 14).Value = -26580

$ws.Cells.Item(140, 8).Value = 7323
$ws.Cells.Item(140, 9).Value = 7323
$ws.Cells.Item(140, 11).Value = 21969
$ws.Cells.Item(140, 13).Value = -16789

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 58513
$ws.Cells.Item(102, 9).Value = 67431.836
$ws.Cells.Item(102, 11).Value = 67431.836
$ws.Cells.Item(102, 13).Value = -65809.836

$ws.Cells.Item(122, 8).Value = 5999.5
$ws.Cells.Item(122, 9).Value = 4833.3335
$ws.Cells.Item(122, 10).Value = 7165.6665
$ws.Cells.Item(122, 11).Value = 14500.0005
$ws.Cells.Item(122, 12).Value = 21496.9995
$ws.Cells.Item(122, 13).Value = -12050.0005
$ws.Cells.Item(122, 14).Value = -26396.9995

$ws.Cells.Item(126, 8).Value = 12855.429
$ws.Cells.Item(126, 9).Value = 9997
$ws.Cells.Item(126, 11).Value = 29991
$ws.Cells.Item(126, 13).Value = -27521

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 8009.3335
$ws.Cells.Item(42, 10).Value = 9514
$ws.Cells.Item(42, 12).Value = 9514
$ws.Cells.Item(42, 14).Value = -10640

$ws.Cells.Item(49, 8).Value = 8009.3335
$ws.Cells.Item(49, 10).Value = 9514
$ws.Cells.Item(49, 12).Value = 9514
$ws.Cells.Item(49, 14).Value = -9808

$ws.Cells.Item(61, 8).Value = 998
$ws.Cells.Item(61, 9).Value = 998
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 998
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -796
$ws.Cells.Item(61, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 998
$ws.Cells.Item(113, 9).Value = 998
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 998
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 1172
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 1000
$ws.Cells.Item(132, 9).Value = 1000
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 3000
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -470
$ws.Cells.Item(132, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(57, 8).Value = 166855.67
$ws.Cells.Item(57, 9).Value = 100177
$ws.Cells.Item(57, 10).Value = 200195
$ws.Cells.Item(57, 11).Value = 100177
$ws.Cells.Item(57, 12).Value = 200195
$ws.Cells.Item(57, 13).Value = -99423
$ws.Cells.Item(57, 14).Value = -201703

$ws.Cells.Item(74, 8).Value = 20999.5
$ws.Cells.Item(74, 10).Value = 21000
$ws.Cells.Item(74, 12).Value = 21000
$ws.Cells.Item(74, 14).Value = -22872

$ws.Cells.Item(77, 8).Value = 20999.5
$ws.Cells.Item(77, 10).Value = 21000
$ws.Cells.Item(77, 12).Value = 21000
$ws.Cells.Item(77, 14).Value = -72360

$ws.Cells.Item(132, 8).Value = 1392.5
$ws.Cells.Item(132, 9).Value = 785
$ws.Cells.Item(132, 11).Value = 2355
$ws.Cells.Item(132, 13).Value = 175

$ws.Cells.Item(136, 8).Value = 11497.546
$ws.Cells.Item(136, 9).Value = 10370.5
$ws.Cells.Item(136, 10).Value = 12141.571
$ws.Cells.Item(136, 11).Value = 31111.5
$ws.Cells.Item(136, 12).Value = 36424.713
$ws.Cells.Item(136, 13).Value = -28561.5
$ws.Cells.Item(136, 14).Value = -41524.713
